# SURE_Poster.pptx edit
# Commit message: "updated AST count in results"
#
# 1) Slide content: the "A total of ### ASTs were generated..." bullet had
#    its highlighted placeholder run ("###", shown in red) filled in with
#    the real figure (153), collapsing the three runs of that sentence back
#    into a single plain run.
# 2) The auto-updating "datetimeFigureOut" date field baked into the slide
#    master and every slide layout advanced from 8/8/2025 to 8/10/2025 (the
#    normal side effect of PowerPoint refreshing that field on save).

$p = $ppt.ActivePresentation

# --- 1. Fix the AST count sentence on the slide(s) ---------------------
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($vi = 1; $vi -le $slide.Shapes.Count; $vi++) {
        $shp = $slide.Shapes.Item($vi)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            $tr = $shp.TextFrame.TextRange
            $full = $tr.Text
            if ($full.Contains("A total of") -and $full.Contains("ASTs were generated")) {
                $startIdx = $full.IndexOf("A total of")
                $len = $full.Length - $startIdx
                # Re-write just this sentence (tail of the paragraph) as one
                # run, replacing the red "###" placeholder with the real
                # count -- matches PowerPoint's own run-merging when you
                # type over a multi-run selection.
                $sub = $tr.Characters($startIdx + 1, $len)
                $sub.Text = "A total of 153 ASTs were generated for use in pattern prediction"
            }
        }
    }
}

# --- 2. Refresh the "datetimeFigureOut" date field everywhere it lives -
$newDate = "8/10/2025"
$ppPlaceholderDate = 16

function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDatePlaceholder = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }
        if (-not $isDatePlaceholder -and $sh.Name -like "Date Placeholder*") {
            $isDatePlaceholder = $true
        }
        if ($isDatePlaceholder) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes

$layouts = $master.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    Update-DatePlaceholders $layouts.Item($L).Shapes
}

Write-Host "Updated AST count and refreshed date fields."
